# Update the workload distribution on the "workload" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Glossy reflections" row (row 25): move 50 points of work from member1 (D)
# over to member3 (F), while leaving member2 (E) untouched.
$ws.Range("D25").Value = 50
$ws.Range("F25").Value = 50

# Restore the cursor/selection to cell E19, matching the author's saved view.
$ws.Range("E19").Select()
